$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 1.930412267966638
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B2:K2").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B3:K3").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 1.995374455241354
$row[0,2] = 1.680409698197509
$row[0,3] = -18.41877499633192
$row[0,4] = 3.448532780533182
$row[0,5] = -18.41877499633192
$row[0,6] = 1.308495697766984
$row[0,7] = -18.41877499633192
$row[0,8] = 1.241264433199964
$row[0,9] = -18.41877499633192
$ws.Range("B4:K4").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 1.673250424238556
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = 2.882824460861115
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B5:K5").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B6:K6").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 2.467273554487778
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B7:K7").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = 1.778375691432487
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B8:K8").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 3.85499016024276
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B9:K9").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = 2.106157862340473
$ws.Range("B10:K10").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = 2.946042024651039
$row[0,4] = -18.41877499633192
$row[0,5] = 2.719874914684135
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = 1.993535377874324
$ws.Range("B11:K11").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B12:K12").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = 2.514495670022804
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = 1.695537984542231
$row[0,9] = 1.883717788996711
$ws.Range("B13:K13").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = 1.558769068337616
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = 1.897087677572775
$ws.Range("B14:K14").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = 1.78503300075582
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B15:K15").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = -18.41877499633192
$row[0,7] = -18.41877499633192
$row[0,8] = 1.413795203277446
$row[0,9] = -18.41877499633192
$ws.Range("B16:K16").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 2.179478369132691
$row[0,2] = 1.877344714874104
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = 2.256524810604284
$row[0,7] = -18.41877499633192
$row[0,8] = 2.693120327430695
$row[0,9] = -18.41877499633192
$ws.Range("B17:K17").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = -18.41877499633192
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = 2.041045813448507
$row[0,7] = -18.41877499633192
$row[0,8] = 2.396555915020221
$row[0,9] = -18.41877499633192
$ws.Range("B18:K18").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = -18.41877499633192
$row[0,2] = 1.975475198958575
$row[0,3] = -18.41877499633192
$row[0,4] = -18.41877499633192
$row[0,5] = -18.41877499633192
$row[0,6] = 1.687354671842746
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B19:K19").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 1.007303397132416
$row[0,2] = 1.483324793126806
$row[0,3] = -18.41877499633192
$row[0,4] = 3.183126792196947
$row[0,5] = -18.41877499633192
$row[0,6] = 1.403125078947964
$row[0,7] = 4.321924183822886
$row[0,8] = -18.41877499633192
$row[0,9] = 2.103516300966614
$ws.Range("B20:K20").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18.41877499633192
$row[0,1] = 1.305534004809138
$row[0,2] = -18.41877499633192
$row[0,3] = 1.65501224205812
$row[0,4] = -18.41877499633192
$row[0,5] = 2.59356105433339
$row[0,6] = 1.466708390770466
$row[0,7] = -18.41877499633192
$row[0,8] = -18.41877499633192
$row[0,9] = -18.41877499633192
$ws.Range("B21:K21").Value = $row
